$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "595.57")
# are stored as literal text instead of being parsed into floating point numbers,
# matching the inlineStr cells in the source workbook. Style is restored to the
# sheet default afterwards so no stray formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.617.36'
$ws.Range("D3").Value = '3.786.35'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '595.57'
$ws.Range("D6").Value = '166.67'
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '3.782.13'
$ws.Range("E7").Value = '  +1.16%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '36.36'
$ws.Range("E14").Value = '  +0.91%  '
$ws.Range("D15").Value = '4.424.56'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").Value = '3.780.51'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").Value = '18.53'
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("D18").Value = '67.598.63'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '7.00'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '10.20'
$ws.Range("E21").Value = '  -4.71%  '
$ws.Range("D22").Value = '457.22'
$ws.Range("E22").Value = '  -1.86%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E24").Value = '  +8.48%  '
$ws.Range("D25").Value = '83.52'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = '11.90'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = '7.29'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '2.20'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '29.83'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").Value = '3.741.48'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").Value = '5.76'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D44").Value = '45.00'
$ws.Range("E44").Value = '  +5.34%  '
$ws.Range("D45").Value = '0.299'
$ws.Range("D46").Value = '47.04'
$ws.Range("E46").Value = '  +2.66%  '
$ws.Range("D47").Value = '8.35'
$ws.Range("E47").Value = '  -2.56%  '
$ws.Range("D48").Value = '148.93'
$ws.Range("E48").Value = '  +2.13%  '
$ws.Range("E49").Value = '  -4.21%  '
$ws.Range("D50").Value = '389.53'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '25.62'
$ws.Range("E51").Value = '  +1.72%  '

# Restore default (General) styling on column D now that the text values are set.
$ws.Range("D2:D51").Style = "Normal"
